$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvoiceLog")

# Row 2 - fill in invoice details (description split into multiple fields:
# date, from, bill to, senders address, receivers address, inv no, company vat reg)
$d = Get-Date -Year 2023 -Month 11 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("B2").Value = $d.Date
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = "V.I.B. Bookshop CC"
$ws.Range("D2").Value = "UNC"
$ws.Range("E2").Value = "U.W.C. Student Centre Modderdam Road, Bellville 7500 PO Box 278, Kasselvlei 7533"
$ws.Range("F2").Value = "ROBERT SORUKNE RD BEUVLUG 7535"
$ws.Range("H2").Value = 8365
$ws.Range("I2").Value = 4100222985
